$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving the cell's
# original style (avoids Excel auto-adding a "quote prefix" style for
# numeric-looking strings, and keeps other cells byte-for-byte stable).
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "275.40"
Set-TextValue "D3" "20.90"
Set-TextValue "D4" "6.192"
Set-TextValue "D5" "0.06183"
Set-TextValue "D6" "3.579"
Set-TextValue "D7" "6.567"
Set-TextValue "D8" "1.499"
Set-TextValue "D9" "0.8246"
Set-TextValue "D10" "0.01380"
Set-TextValue "D11" "0.1615"
Set-TextValue "D12" "0.08215"
Set-TextValue "D13" "0.03535"
Set-TextValue "D14" "0.03104"
Set-TextValue "D15" "0.09130"
Set-TextValue "D16" "3.716"
Set-TextValue "D17" "0.001607"
Set-TextValue "D18" "0.04699"
Set-TextValue "D19" "0.006460"
Set-TextValue "D21" "0.0001501"
Set-TextValue "D22" "3.803"
Set-TextValue "D23" "2.282"
Set-TextValue "D28" "0.0001581"
Set-TextValue "E28" "27UpBotsUBXT"
Set-TextValue "D40" "0.04631"
Set-TextValue "D41" "0.007025"
Set-TextValue "E41" "40KickTokenKICKBestin24h"
Set-TextValue "D42" "0.004604"
Set-TextValue "D43" "0.1100"
Set-TextValue "D44" "0.01084"
Set-TextValue "D45" "0.00006162"
Set-TextValue "D47" "0.8459"
Set-TextValue "D48" "0.002571"
Set-TextValue "D49" "0.00001901"
Set-TextValue "D50" "0.01241"
